# Applies the "More changes made to the data_set" commit:
#  - Removes the stray "MaterialCost" column (P) entirely
#  - Fills in previously-blank Details/Subsystem cells for rows 22-24
#  - Inserts a new "Battery Management System" line item as row 26
#    (pushing the old Sensor/DAQ row down to row 27)
#  - Adjusts column widths (D, O) and the active-cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the erroneous "MaterialCost" column P entirely (header only had data in P1).
$ws.Columns.Item(16).Delete()

# 2. Backfill Details (M) / Subsystem (O) for the Electric Motor / Electric Cable / Battery rows
#    which previously had some blanks.
$ws.Range("M22").Value = "Supplier Quote"
$ws.Range("O22").Value = "Tractive System"
$ws.Range("M23").Value = "Supplier Quote"
$ws.Range("O23").Value = "Tractive System"
$ws.Range("M24").Value = "Supplier Quote"
$ws.Range("O24").Value = "Tractive System"

# 3. Move the existing row 26 (Sensor / DAQ) down to row 27 so a new row can be
#    inserted in its place, then re-populate row 27 with that original data
#    (the formula is re-written relative to its new row).
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "DAQ"
$ws.Range("C27").Value = "ABC123-PQR"
$ws.Range("D27").Value = "Sensor"
$ws.Range("E27").Value = "Engine Speed Sensor"
$ws.Range("F27").Value = 100
$ws.Range("G27").Value = 30
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 30
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = 1
$ws.Range("L27").Formula = "=F27*G27+H27+I27+J27+K27"
$ws.Range("M27").Value = "Among Mens"
$ws.Range("N27").Value = "Datasheet"
$ws.Range("O27").Value = "DAQ"

# 4. Overwrite row 26 in place with the new "Battery Management System" line item.
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "Tractive System"
$ws.Range("C26").Value = "DEF345-GHI"
$ws.Range("D26").Value = "Battery Management System"
$ws.Range("E26").Value = "Orion BMS 2"
$ws.Range("F26").Value = 1230
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 30
$ws.Range("I26").Value = 30
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 60
$ws.Range("L26").Formula = "=F26*G26+H26+I26+J26+K26"
$ws.Range("M26").Value = "Critical Component"
$ws.Range("O26").Value = "Tractive System"

# 5. Column width tweaks: widen column D, and give column O an explicit width
#    (closest attainable values - the host quantizes ColumnWidth to 1/6-character steps).
$ws.Columns.Item(4).ColumnWidth = 25.333333333333332
$ws.Columns.Item(15).ColumnWidth = 15.666666666666666

# 6. Restore the selection to where the editor last left off.
$ws.Range("N26").Select()
